$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full ifoCAST series evaluation: each row's B:F values shift up by one
# (row N now holds what used to be row N+1), N (col G) decreases by one,
# and the final row (11) receives the newly-added evaluation figures.

$data = @(
    @(-0.01572770659491909, 0.4445496780545906, 0.3664663352916562, 0.6053646300302457, 0.6280045102432898, 14),
    @(-0.06939138655364743, 0.5952380793966869, 0.6460436807749044, 0.803768424843191,  0.8334651917952676, 13),
    @(0.06184427096005268,  0.33700112578256,   0.209417853196624,  0.4576219544521701, 0.4735857100263252, 12),
    @(-0.06798794410060216, 0.4968139668552824, 0.377721656252307,  0.6145906412013666, 0.640631895303542,  11),
    @(0.01694414359772521,  0.2607354299704002, 0.1060912636732037, 0.3257165388389169, 0.3428704965612517, 10),
    @(-0.06054918682133375, 0.3330393258437492, 0.1559239451819853, 0.3948720617896198, 0.4138719111052573, 9),
    @(-0.01715367655228583, 0.2997021343775136, 0.1700925356682539, 0.4124227632760514, 0.4405169522264152, 8),
    @(-0.05536040423422233, 0.4661046776547557, 0.3164365824987914, 0.5625269615749909, 0.6046490245027418, 7),
    @(-0.0196809274911882,  0.3184982995745606, 0.1151557242996088, 0.3393460244346599, 0.371109232533003,  6),
    @(0.01522744611739815,  0.4477926050010403, 0.2261067416802728, 0.475506826113225,  0.5313601257210295, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
